# Apply updated classification-report values (DecisionTreeClassifier / lgbm TPE search results)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B7" = 0.5079365079365079
    "C7" = 0.8421052631578947
    "D7" = 0.6336633663366337

    "B8" = 0.8
    "C8" = 0.4363636363636363
    "D8" = 0.5647058823529413

    "B9" = 0.6021505376344086
    "C9" = 0.6021505376344086
    "D9" = 0.6021505376344086
    "E9" = 0.6021505376344086

    "B10" = 0.653968253968254
    "C10" = 0.6392344497607655
    "D10" = 0.5991846243447875

    "B11" = 0.6806622290493258
    "C11" = 0.6021505376344086
    "D11" = 0.5928820586043425

    "B12" = 0.5609756097560976
    "C12" = 0.6052631578947368
    "D12" = 0.5822784810126583

    "B13" = 0.7115384615384616
    "C13" = 0.6727272727272727
    "D13" = 0.6915887850467289

    "B14" = 0.6451612903225806
    "C14" = 0.6451612903225806
    "D14" = 0.6451612903225806
    "E14" = 0.6451612903225806

    "B15" = 0.6362570356472796
    "C15" = 0.6389952153110048
    "D15" = 0.6369336330296936

    "B16" = 0.6500181565091085
    "C16" = 0.6451612903225806
    "D16" = 0.6469243597424851

    "B22" = 0.8
    "C22" = 0.631578947368421
    "D22" = 0.7058823529411765

    "B23" = 0.7777777777777778
    "C23" = 0.8909090909090909
    "D23" = 0.8305084745762712

    "B24" = 0.7849462365591398
    "C24" = 0.7849462365591398
    "D24" = 0.7849462365591398
    "E24" = 0.7849462365591398

    "B25" = 0.788888888888889
    "C25" = 0.7612440191387559
    "D25" = 0.7681954137587239

    "B26" = 0.7868578255675029
    "C26" = 0.7849462365591398
    "D26" = 0.7795859732630066
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
